# NPC.xlsx — "add drop item list record"
#
# Record_PosList sheet (sheet 2): row 1 headers L1:O1 were mistakenly all
# "float" and the real header labels (X / Y / Z / StayTime) were sitting
# in the example data row (L2:O2) instead. This swaps them so row 1 holds
# the column headers and row 2 holds an example "float" data row — and
# moves/retexts the matching cell comments accordingly. Also flips which
# sheet/cell is active/selected.

$wb = $excel.ActiveWorkbook
$wsProperty = $wb.Worksheets.Item("Property")
$wsPosList  = $wb.Worksheets.Item("Record_PosList")

# --- Record_PosList!L1:O1 <-> L2:O2 -----------------------------------
# Row 1 becomes the real header labels, row 2 becomes the "float" sample row.
$wsPosList.Range("L1").Value = "X"
$wsPosList.Range("M1").Value = "Y"
$wsPosList.Range("N1").Value = "Z"
$wsPosList.Range("O1").Value = "StayTime"

$wsPosList.Range("L2").Value = "float"
$wsPosList.Range("M2").Value = "float"
$wsPosList.Range("N2").Value = "float"
$wsPosList.Range("O2").Value = "float"

# --- Move the 4 cell comments from row 1 to row 2, with updated text --
$wsPosList.Range("L1").Comment.Delete()
$wsPosList.Range("M1").Comment.Delete()
$wsPosList.Range("N1").Comment.Delete()
$wsPosList.Range("O1").Comment.Delete()

$wsPosList.Range("L2").AddComment("强化等级")
$wsPosList.Range("M2").AddComment("强化等级")
$wsPosList.Range("N2").AddComment("强化等级")
$wsPosList.Range("O2").AddComment("镶嵌宝石，逗号分隔")

# --- Selection / active-sheet swap ------------------------------------
# Property used to be the tab in front with A40 selected; now
# Record_PosList is in front (with O8 selected) and Property keeps a
# J33 selection in the background.
$wsProperty.Range("J33").Select()
$wsPosList.Activate()
$wsPosList.Range("O8").Select()
